$d = $word.ActiveDocument

# 1. Replace the italicized magazine/blog title text.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("the 221B Blog", $false, $false, $false, $false, $false, $true, 1, $false, "Twin Life Magazine", 2)

# 2. Replace the page range "927 - 929" with "923 - 924"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("927", $false, $false, $false, $false, $false, $true, 1, $false, "923", 2)

$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$find3.Execute("929", $false, $false, $false, $false, $false, $true, 1, $false, "924", 2)

Write-Output "done"
